$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 130, pushing existing rows 130-214 down to 132-216
$ws.Rows("130:131").Insert()

# Fill the first newly inserted row (130) with the new weekly price entry
$ws.Range("A130").Value = 3
$ws.Range("B130").Value = "Femacal de La Calera"
$ws.Range("C130").Value = "Coquimbo"
$ws.Range("D130").Value = 44488
$ws.Range("E130").Value = 5
$ws.Range("F130").Value = 100112013
$ws.Range("G130").Value = "Alcachofa"
$ws.Range("H130").Value = "Española"
$ws.Range("I130").Value = "Extra"
$ws.Range("J130").Value = 3500
$ws.Range("K130").Value = 350
$ws.Range("L130").Value = 350
$ws.Range("M130").Value = 350
$ws.Range("N130").Value = "$/unidad"
$ws.Range("O130").Value = "Llay Llay"
$ws.Range("P130").Value = 350
$ws.Range("Q130").Value = 1
$ws.Range("R130").Value = "Hortaliza"

# Fill the second newly inserted row (131) with the new weekly price entry
$ws.Range("A131").Value = 3
$ws.Range("B131").Value = "Femacal de La Calera"
$ws.Range("C131").Value = "Coquimbo"
$ws.Range("D131").Value = 44488
$ws.Range("E131").Value = 5
$ws.Range("F131").Value = 100112013
$ws.Range("G131").Value = "Alcachofa"
$ws.Range("H131").Value = "Española"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 3800
$ws.Range("K131").Value = 300
$ws.Range("L131").Value = 300
$ws.Range("M131").Value = 300
$ws.Range("N131").Value = "$/unidad"
$ws.Range("O131").Value = "Llay Llay"
$ws.Range("P131").Value = 300
$ws.Range("Q131").Value = 1
$ws.Range("R131").Value = "Hortaliza"
